# [PV-94][WIP] Support for plans without sticky-ids or levels
# Rename column headers on the plan import sheet to match the new
# "no sticky-id / no level" header scheme, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

$ws.Activate()
$ws.Range("F1").Select()
